$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) for columns M, N, O ---
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Match the formatting already used by the other header cells (bold, bordered,
# centered) without inventing a brand new style entry - copy format from L1.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Populate the new columns for each data row (2-19) ---
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 13).Value = "before"
    $ws.Cells.Item($r, 14).Value = 20160636
    $ws.Cells.Item($r, 15).Value = 3
}
